$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-analysis correction: the "% of DMSO control" column (C) originally
# chained each row off the previous row's raw absorbance in column B
# (B_n / B_(n-1) * C_(n-1)), which is wrong - it should always divide by
# the DMSO control absorbance in B48 (0.8648) to get a true percentage of
# control. Fix C49 and propagate the corrected formula down through the
# C50:C57 shared-formula block.
$ws.Range("C49").Formula = "=B49/0.8648*C48"
$ws.Range("C50:C57").Formula = "=B50/0.8648*C49"

# Daily entry: leave the selection where work stopped for the day.
$ws.Range("E54").Select()
